$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 169; existing rows 169..259 shift down to 170..260
$ws.Rows(169).Insert()

# Populate the newly inserted row 169 with the new weekly data point
$ws.Cells.Item(169, 1).Value = 10
$ws.Cells.Item(169, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(169, 3).Value = "La Araucanía"
$ws.Cells.Item(169, 4).Value = 45016
$ws.Cells.Item(169, 5).Value = 9
$ws.Cells.Item(169, 6).Value = 100112012
$ws.Cells.Item(169, 7).Value = "Espinaca"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 20
$ws.Cells.Item(169, 11).Value = 10000
$ws.Cells.Item(169, 12).Value = 10000
$ws.Cells.Item(169, 13).Value = 10000
$ws.Cells.Item(169, 14).Value = "`$/docena de atados"
$ws.Cells.Item(169, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(169, 16).Value = 3333
$ws.Cells.Item(169, 17).Value = 3
$ws.Cells.Item(169, 18).Value = "Hortaliza"
